# Update NATMI LR-pair TPM results (Ccl11 -> Ackr4) with the new TPM-derived
# scores and expand the table from the partial 12-row (missing the
# "Resolving-Mac" target / excluding same-cluster target rows) matrix to the
# full 4x4 Sending-cluster x Target-cluster matrix (16 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "ECs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.738254
$arr[0,7] = 2.214762
$arr[0,8] = 0.005691320045803731
$arr[0,9] = 0.005691320045803731
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.1143813333333333
$arr[0,13] = 0.343144
$arr[0,14] = 0.128300337591142
$arr[0,15] = 0.1283003375911419
$arr[0,16] = 0.08444247685866667
$arr[0,17] = 0.759982291728
$arr[0,18] = 0.0007301982832158522
$arr[0,19] = 0.0007301982832158521
$ws.Range("A2:T2").Value = $arr

# Row 3: ECs -> FAPs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "FAPs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.738254
$arr[0,7] = 2.214762
$arr[0,8] = 0.005691320045803731
$arr[0,9] = 0.005691320045803731
$arr[0,10] = 3.0
$arr[0,11] = 1.0
$arr[0,12] = 0.7200953333333334
$arr[0,13] = 2.160286
$arr[0,14] = 0.8077233554817153
$arr[0,15] = 0.8077233554817151
$arr[0,16] = 0.5316132602146667
$arr[0,17] = 4.784519341932
$arr[0,18] = 0.004597012124516939
$arr[0,19] = 0.004597012124516939
$ws.Range("A3:T3").Value = $arr

# Row 4: ECs -> MuSCs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.738254
$arr[0,7] = 2.214762
$arr[0,8] = 0.005691320045803731
$arr[0,9] = 0.005691320045803731
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.03357866666666667
$arr[0,13] = 0.100736
$arr[0,14] = 0.03766483694187069
$arr[0,15] = 0.03766483694187069
$arr[0,16] = 0.02478958498133333
$arr[0,17] = 0.223106264832
$arr[0,18] = 0.0002143626415091975
$arr[0,19] = 0.0002143626415091975
$ws.Range("A4:T4").Value = $arr

# Row 5: ECs -> Resolving-Mac (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "ECs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.738254
$arr[0,7] = 2.214762
$arr[0,8] = 0.005691320045803731
$arr[0,9] = 0.005691320045803731
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.023457
$arr[0,13] = 0.070371
$arr[0,14] = 0.02631146998527222
$arr[0,15] = 0.02631146998527222
$arr[0,16] = 0.017317224078
$arr[0,17] = 0.155855016702
$arr[0,18] = 0.000149746996561743
$arr[0,19] = 0.000149746996561743
$ws.Range("A5:T5").Value = $arr

# Row 6: FAPs -> ECs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "ECs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 122.458089
$arr[0,7] = 367.374267
$arr[0,8] = 0.9440493064670392
$arr[0,9] = 0.9440493064670391
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.1143813333333333
$arr[0,13] = 0.343144
$arr[0,14] = 0.128300337591142
$arr[0,15] = 0.1283003375911419
$arr[0,16] = 14.006919497272
$arr[0,17] = 126.062275475448
$arr[0,18] = 0.1211218447224046
$arr[0,19] = 0.1211218447224045
$ws.Range("A6:T6").Value = $arr

# Row 7: FAPs -> FAPs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "FAPs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 122.458089
$arr[0,7] = 367.374267
$arr[0,8] = 0.9440493064670392
$arr[0,9] = 0.9440493064670391
$arr[0,10] = 3.0
$arr[0,11] = 1.0
$arr[0,12] = 0.7200953333333334
$arr[0,13] = 2.160286
$arr[0,14] = 0.8077233554817153
$arr[0,15] = 0.8077233554817151
$arr[0,16] = 88.18149841781802
$arr[0,17] = 793.6334857603621
$arr[0,18] = 0.762530673559743
$arr[0,19] = 0.7625306735597429
$ws.Range("A7:T7").Value = $arr

# Row 8: FAPs -> MuSCs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 122.458089
$arr[0,7] = 367.374267
$arr[0,8] = 0.9440493064670392
$arr[0,9] = 0.9440493064670391
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.03357866666666667
$arr[0,13] = 0.100736
$arr[0,14] = 0.03766483694187069
$arr[0,15] = 0.03766483694187069
$arr[0,16] = 4.111979351168
$arr[0,17] = 37.00781416051201
$arr[0,18] = 0.03555746319316714
$arr[0,19] = 0.03555746319316713
$ws.Range("A8:T8").Value = $arr

# Row 9: FAPs -> Resolving-Mac (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 122.458089
$arr[0,7] = 367.374267
$arr[0,8] = 0.9440493064670392
$arr[0,9] = 0.9440493064670391
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.023457
$arr[0,13] = 0.070371
$arr[0,14] = 0.02631146998527222
$arr[0,15] = 0.02631146998527222
$arr[0,16] = 2.872499393673001
$arr[0,17] = 25.852494543057
$arr[0,18] = 0.02483932499172456
$arr[0,19] = 0.02483932499172455
$ws.Range("A9:T9").Value = $arr

# Row 10: MuSCs -> ECs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "ECs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 5.698467
$arr[0,7] = 17.095401
$arr[0,8] = 0.0439304080539368
$arr[0,9] = 0.04393040805393679
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.1143813333333333
$arr[0,13] = 0.343144
$arr[0,14] = 0.128300337591142
$arr[0,15] = 0.1283003375911419
$arr[0,16] = 0.651798253416
$arr[0,17] = 5.866184280744
$arr[0,18] = 0.005636286183836712
$arr[0,19] = 0.005636286183836711
$ws.Range("A10:T10").Value = $arr

# Row 11: MuSCs -> FAPs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "FAPs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 5.698467
$arr[0,7] = 17.095401
$arr[0,8] = 0.0439304080539368
$arr[0,9] = 0.04393040805393679
$arr[0,10] = 3.0
$arr[0,11] = 1.0
$arr[0,12] = 0.7200953333333334
$arr[0,13] = 2.160286
$arr[0,14] = 0.8077233554817153
$arr[0,15] = 0.8077233554817151
$arr[0,16] = 4.103439493854
$arr[0,17] = 36.930955444686
$arr[0,18] = 0.0354836166010068
$arr[0,19] = 0.03548361660100679
$ws.Range("A11:T11").Value = $arr

# Row 12: MuSCs -> MuSCs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 5.698467
$arr[0,7] = 17.095401
$arr[0,8] = 0.0439304080539368
$arr[0,9] = 0.04393040805393679
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.03357866666666667
$arr[0,13] = 0.100736
$arr[0,14] = 0.03766483694187069
$arr[0,15] = 0.03766483694187069
$arr[0,16] = 0.191346923904
$arr[0,17] = 1.722122315136
$arr[0,18] = 0.001654631656141372
$arr[0,19] = 0.001654631656141372
$ws.Range("A12:T12").Value = $arr

# Row 13: MuSCs -> Resolving-Mac (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "MuSCs"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 5.698467
$arr[0,7] = 17.095401
$arr[0,8] = 0.0439304080539368
$arr[0,9] = 0.04393040805393679
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.023457
$arr[0,13] = 0.070371
$arr[0,14] = 0.02631146998527222
$arr[0,15] = 0.02631146998527222
$arr[0,16] = 0.133668940419
$arr[0,17] = 1.203020463771
$arr[0,18] = 0.001155873612951919
$arr[0,19] = 0.001155873612951919
$ws.Range("A13:T13").Value = $arr

# Row 14: Resolving-Mac -> ECs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "ECs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.8209666666666666
$arr[0,7] = 2.4629
$arr[0,8] = 0.006328965433220369
$arr[0,9] = 0.006328965433220369
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.1143813333333333
$arr[0,13] = 0.343144
$arr[0,14] = 0.128300337591142
$arr[0,15] = 0.1283003375911419
$arr[0,16] = 0.09390326195555555
$arr[0,17] = 0.8451293576
$arr[0,18] = 0.0008120084016848413
$arr[0,19] = 0.0008120084016848411
$ws.Range("A14:T14").Value = $arr

# Row 15: Resolving-Mac -> FAPs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "FAPs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.8209666666666666
$arr[0,7] = 2.4629
$arr[0,8] = 0.006328965433220369
$arr[0,9] = 0.006328965433220369
$arr[0,10] = 3.0
$arr[0,11] = 1.0
$arr[0,12] = 0.7200953333333334
$arr[0,13] = 2.160286
$arr[0,14] = 0.8077233554817153
$arr[0,15] = 0.8077233554817151
$arr[0,16] = 0.591174265488889
$arr[0,17] = 5.3205683894
$arr[0,18] = 0.005112053196448544
$arr[0,19] = 0.005112053196448543
$ws.Range("A15:T15").Value = $arr

# Row 16: Resolving-Mac -> MuSCs (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "MuSCs"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.8209666666666666
$arr[0,7] = 2.4629
$arr[0,8] = 0.006328965433220369
$arr[0,9] = 0.006328965433220369
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.03357866666666667
$arr[0,13] = 0.100736
$arr[0,14] = 0.03766483694187069
$arr[0,15] = 0.03766483694187069
$arr[0,16] = 0.02756696604444444
$arr[0,17] = 0.2481026944
$arr[0,18] = 0.0002383794510529812
$arr[0,19] = 0.0002383794510529812
$ws.Range("A16:T16").Value = $arr

# Row 17: Resolving-Mac -> Resolving-Mac (Ccl11 -> Ackr4)
$arr = New-Object 'object[,]' 1,20
$arr[0,0] = "Resolving-Mac"
$arr[0,1] = "Ccl11"
$arr[0,2] = "Ackr4"
$arr[0,3] = "Resolving-Mac"
$arr[0,4] = 3.0
$arr[0,5] = 1.0
$arr[0,6] = 0.8209666666666666
$arr[0,7] = 2.4629
$arr[0,8] = 0.006328965433220369
$arr[0,9] = 0.006328965433220369
$arr[0,10] = 1.0
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.023457
$arr[0,13] = 0.070371
$arr[0,14] = 0.02631146998527222
$arr[0,15] = 0.02631146998527222
$arr[0,16] = 0.0192574151
$arr[0,17] = 0.1733167359
$arr[0,18] = 0.0001665243840340031
$arr[0,19] = 0.0001665243840340031
$ws.Range("A17:T17").Value = $arr
